$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '68.498.75'
$ws.Range('E2').Value = '  -1.62%  '
$ws.Range('D3').Value = '2.458.25'
$ws.Range('E3').Value = '  -1.76%  '
$ws.Range('E4').Value = '  -0.04%  '
Set-TextValue $ws 'D5' '560.12'
$ws.Range('E5').Value = '  -2.51%  '
Set-TextValue $ws 'D6' '164.30'
$ws.Range('E6').Value = '  -1.24%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.57%  '
$ws.Range('D9').Value = '2.457.33'
$ws.Range('E9').Value = '  -1.76%  '
$ws.Range('E10').Value = '  -6.09%  '
$ws.Range('E11').Value = '  -1.81%  '
$ws.Range('E12').Value = '  -5.62%  '
Set-TextValue $ws 'D13' '4.83'
$ws.Range('E13').Value = '  -2.12%  '
$ws.Range('D14').Value = '2.904.82'
$ws.Range('E14').Value = '  -1.81%  '
$ws.Range('D15').Value = '68.355.56'
$ws.Range('E15').Value = '  -1.60%  '
$ws.Range('E16').Value = '  -3.66%  '
Set-TextValue $ws 'D17' '23.35'
$ws.Range('E17').Value = '  -5.39%  '
$ws.Range('D18').Value = '2.513.48'
$ws.Range('E18').Value = '  +0.40%  '
Set-TextValue $ws 'D19' '11.01'
$ws.Range('E19').Value = '  -1.75%  '
Set-TextValue $ws 'D20' '344.87'
Set-TextValue $ws 'D21' '7.20'
$ws.Range('E21').Value = '  -4.16%  '
Set-TextValue $ws 'D22' '3.78'
$ws.Range('E22').Value = '  -2.97%  '
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('E24').Value = '  -3.28%  '
Set-TextValue $ws 'D25' '67.87'
$ws.Range('E25').Value = '  -3.95%  '
$ws.Range('E26').Value = '  +9.44%  '
Set-TextValue $ws 'D27' '3.72'
$ws.Range('E27').Value = '  -5.18%  '
$ws.Range('D28').Value = '2.586.53'
$ws.Range('E28').Value = '  -1.65%  '
$ws.Range('E29').Value = '  -6.65%  '
$ws.Range('D30').Value = '0.0₃0836'
$ws.Range('E30').Value = '  -6.07%  '
Set-TextValue $ws 'D31' '7.22'
$ws.Range('E31').Value = '  -8.08%  '
Set-TextValue $ws 'D32' '3.44'
$ws.Range('E32').Value = '  +136.61%  '
Set-TextValue $ws 'D33' '434.37'
$ws.Range('E33').Value = '  -5.09%  '
Set-TextValue $ws 'D34' '0.999'
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -3.08%  '
Set-TextValue $ws 'D36' '1.66'
$ws.Range('E36').Value = '  -3.94%  '
Set-TextValue $ws 'D37' '157.87'
$ws.Range('E37').Value = '  +0.36%  '
$ws.Range('E38').Value = '  -0.22%  '
Set-TextValue $ws 'D39' '0.110'
$ws.Range('E39').Value = '  -4.69%  '
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('E41').Value = '  -2.37%  '
Set-TextValue $ws 'D42' '0.306'
$ws.Range('E42').Value = '  -3.35%  '
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('E44').Value = '  -4.84%  '
Set-TextValue $ws 'D45' '1.10'
$ws.Range('E45').Value = '  +1.96%  '
Set-TextValue $ws 'D47' '134.95'
$ws.Range('E47').Value = '  -4.30%  '
Set-TextValue $ws 'D48' '3.36'
$ws.Range('E48').Value = '  -3.57%  '
$ws.Range('E49').Value = '  -1.89%  '
$ws.Range('E50').Value = '  -6.55%  '
$ws.Range('E51').Value = '  -2.40%  '

Write-Host "Applied 77 cell updates"
